$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates derived from the crypto price refresh diff.
# Cells in $textCells hold numeric-looking strings (e.g. '0.500') that must
# stay text (matching the workbook's original inlineStr storage) instead of
# being auto-coerced to numbers by Excel (which would also drop trailing zeros).
$textCells = @("D5", "D6", "D10", "D11", "D13", "D15", "D17", "D19", "D23", "D25", "D29", "D30", "D33", "D38", "D42", "D43", "D45", "D46", "D47", "D49", "D50", "D51")

$updates = [ordered]@{
    "D2" = "26.673.50"
    "E2" = "  +1.22%  "
    "D3" = "1.631.23"
    "E3" = "  +0.49%  "
    "E4" = "  +0.02%  "
    "D5" = "213.35"
    "E5" = "  +0.59%  "
    "D6" = "0.500"
    "E6" = "  +3.29%  "
    "E7" = "  +0.01%  "
    "E8" = "  +1.75%  "
    "E9" = "  +0.95%  "
    "D10" = "19.25"
    "E10" = "  +2.44%  "
    "D11" = "0.0841"
    "E11" = "  +3.22%  "
    "D12" = "1.859.13"
    "E12" = "  +0.56%  "
    "B13" = "Polkadot"
    "C13" = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
    "D13" = "4.09"
    "E13" = "  +1.48%  "
    "B14" = "WrappedEther"
    "C14" = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
    "D14" = "1.584.67"
    "E14" = "  -2.30%  "
    "D15" = "0.524"
    "E15" = "  +0.97%  "
    "D16" = "26.650.47"
    "E16" = "  +1.13%  "
    "D17" = "63.50"
    "E17" = "  +1.57%  "
    "E18" = "  +2.36%  "
    "D19" = "218.79"
    "E19" = "  +7.79%  "
    "E21" = "  +0.73%  "
    "E22" = "  +2.06%  "
    "D23" = "9.37"
    "E23" = "  +0.62%  "
    "E24" = "  +4.65%  "
    "D25" = "147.65"
    "E25" = "  +2.30%  "
    "E27" = "  +1.31%  "
    "E28" = "  +4.09%  "
    "D29" = "15.54"
    "E29" = "  +2.09%  "
    "D30" = "0.0504"
    "E30" = "  -3.02%  "
    "E31" = "  +0.35%  "
    "E32" = "  +3.70%  "
    "D33" = "2.98"
    "E33" = "  +2.25%  "
    "E34" = "  +0.88%  "
    "E35" = "  +0.33%  "
    "D36" = "1.214.50"
    "E36" = "  +4.92%  "
    "E37" = "  +4.39%  "
    "D38" = "0.805"
    "E38" = "  -0.56%  "
    "E39" = "  -0.01%  "
    "E40" = "  +0.59%  "
    "E41" = "  -1.87%  "
    "D42" = "0.794"
    "E42" = "  +1.15%  "
    "D43" = "5.36"
    "E43" = "  -0.44%  "
    "D44" = "1.768.13"
    "E44" = "  +0.44%  "
    "D45" = "92.78"
    "E45" = "  +0.20%  "
    "D46" = "1.56"
    "E46" = "  +2.87%  "
    "D47" = "55.02"
    "E47" = "  +2.16%  "
    "B48" = "BabyDogeCoin"
    "C48" = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
    "D48" = "0.0₆0104"
    "E48" = "  -0.56%  "
    "B49" = "Cronos"
    "C49" = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
    "D49" = "0.0511"
    "E49" = "  +0.38%  "
    "B50" = "EnergySwap"
    "C50" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "D50" = "7.62"
    "E50" = "  +3.95%  "
    "B51" = "Mantle"
    "C51" = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
    "D51" = "0.409"
    "E51" = "  -0.01%  "
}

foreach ($cellRef in $updates.Keys) {
    $value = $updates[$cellRef]
    $range = $ws.Range($cellRef)
    if ($textCells -contains $cellRef) {
        # Force text storage, write the value, then drop the temporary
        # number-format override so the cell's style index is unchanged.
        $range.NumberFormat = "@"
        $range.Value = $value
        $range.ClearFormats()
    } else {
        $range.Value = $value
    }
}
